$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (Price and Volume(1h) columns).
# NumberFormat is forced to Text ("@") before the assignment and the style is
# reset back to Normal afterwards so that values such as "131.10" or "0.999"
# are kept as literal strings instead of being auto-coerced into numbers by Excel.
$updates = [ordered]@{
    "D2" = "64.113.96"
    "E2" = "  -0.22%  "
    "D3" = "3.473.01"
    "E3" = "  -0.73%  "
    "D4" = "0.999"
    "E4" = "  -0.07%  "
    "D5" = "583.87"
    "E5" = "  -0.39%  "
    "D6" = "131.10"
    "E6" = "  -2.21%  "
    "E7" = "  +0.04%  "
    "D8" = "0.482"
    "E8" = "  -1.03%  "
    "D9" = "7.64"
    "E9" = "  +5.54%  "
    "E10" = "  -1.21%  "
    "E11" = "  +0.38%  "
    "D12" = "4.062.28"
    "E12" = "  -0.85%  "
    "E13" = "  -0.16%  "
    "E14" = "  -2.37%  "
    "D15" = "3.471.39"
    "E15" = "  -0.85%  "
    "D16" = "64.067.94"
    "E16" = "  -0.32%  "
    "D17" = "24.36"
    "E17" = "  -5.87%  "
    "E18" = "  +0.62%  "
    "E19" = "  -1.36%  "
    "D20" = "13.42"
    "E20" = "  -1.44%  "
    "D21" = "384.10"
    "E21" = "  -2.48%  "
    "E22" = "  -0.66%  "
    "D23" = "3.612.32"
    "E23" = "  -0.74%  "
    "D24" = "74.83"
    "E24" = "  +0.73%  "
    "E25" = "  +0.17%  "
    "E26" = "  +0.41%  "
    "E27" = "  -2.86%  "
    "E28" = "  -0.08%  "
    "E29" = "  -0.54%  "
    "E30" = "  -4.65%  "
    "E31" = "  -4.79%  "
    "E32" = "  -4.12%  "
    "D33" = "0.153"
    "E33" = "  +1.39%  "
    "D34" = "3.500.00"
    "E34" = "  -0.53%  "
    "E35" = "  -0.05%  "
    "D36" = "22.91"
    "E36" = "  -2.20%  "
    "D37" = "5.19"
    "E37" = "  +0.84%  "
    "E38" = "  -2.31%  "
    "E39" = "  -3.72%  "
    "E40" = "  -2.20%  "
    "D41" = "0.0778"
    "E41" = "  -0.65%  "
    "E42" = "  -1.16%  "
    "E43" = "  -0.06%  "
    "D44" = "41.22"
    "E44" = "  -1.61%  "
    "E45" = "  -2.62%  "
    "D46" = "1.62"
    "E46" = "  -2.14%  "
    "D47" = "23.53"
    "E47" = "  -6.78%  "
    "E48" = "  -5.02%  "
    "E49" = "  -0.96%  "
    "D50" = "0.901"
    "E50" = "  +0.61%  "
    "D51" = "2.327.53"
    "E51" = "  -5.29%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
